$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "1.005", "0.3818") are stored verbatim instead of being
# auto-parsed into floating point numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '27.092.62'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '1.781.36'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '336.75'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.3818'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = '0.3408'
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('D9').Value = '47.87'
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').Value = '1.184'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('D11').Value = '0.07439'
$ws.Range('E11').Value = '  -3.56%  '
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').Value = '21.56'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('D14').Value = '6.417'
$ws.Range('E14').Value = '  -3.13%  '
$ws.Range('D15').Value = '1.782.16'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').Value = '7.095'
$ws.Range('E16').Value = '  -1.72%  '
$ws.Range('D17').Value = '0.00001091'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '0.06647'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('D19').Value = '83.31'
$ws.Range('E19').Value = '  -3.17%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').Value = '6.514'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = '17.32'
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').Value = '27.101.29'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').Value = '12.21'
$ws.Range('E24').Value = '  -7.74%  '
$ws.Range('D25').Value = '2.376'
$ws.Range('E25').Value = '  -3.81%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '21.08'
$ws.Range('E26').Value = '  -3.93%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.489'
$ws.Range('E27').Value = '  -6.82%  '
$ws.Range('D28').Value = '1.446'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').Value = '153.90'
$ws.Range('E29').Value = '  -0.41%  '
$ws.Range('D30').Value = '1.983.68'
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('D31').Value = '133.76'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').Value = '3.981'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').Value = '6.003'
$ws.Range('E33').Value = '  -4.97%  '
$ws.Range('D34').Value = '0.08649'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').Value = '12.96'
$ws.Range('E35').Value = '  -6.81%  '
$ws.Range('D36').Value = '1.624'
$ws.Range('E36').Value = '  -4.21%  '
$ws.Range('D37').Value = '5.372'
$ws.Range('E37').Value = '  -4.40%  '
$ws.Range('D38').Value = '0.6803'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('D39').Value = '0.06309'
$ws.Range('E39').Value = '  -2.51%  '
$ws.Range('D40').Value = '0.02321'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').Value = '0.2166'
$ws.Range('E41').Value = '  -4.54%  '
$ws.Range('D42').Value = '1.245'
$ws.Range('E42').Value = '  -4.34%  '
$ws.Range('D43').Value = '8.413'
$ws.Range('E43').Value = '  -5.92%  '
$ws.Range('E44').Value = '  -3.96%  '
$ws.Range('D45').Value = '1.001'
$ws.Range('D46').Value = '0.6382'
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('E47').Value = '  -4.62%  '
$ws.Range('D48').Value = '2.135'
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('D49').Value = '130.90'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('D50').Value = '0.07088'
$ws.Range('E50').Value = '  -3.27%  '
$ws.Range('D51').Value = '78.40'
$ws.Range('E51').Value = '  -2.49%  '

# Restore the default cell style now that the text is committed,
# so no stray style survives the edit (matches original formatting).
$ws.Range('D2:E51').Style = 'Normal'
